$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2315597.2
$ws.Range("I18").Value = 2315597.2
$ws.Range("K18").Value = 2315597.2
$ws.Range("M18").Value = -2315313.2
$ws.Range("H33").Value = 563.5714
$ws.Range("I33").Value = 481.0625
$ws.Range("J33").Value = 673.5833
$ws.Range("K33").Value = 481.0625
$ws.Range("L33").Value = 673.5833
$ws.Range("M33").Value = -252.0625
$ws.Range("N33").Value = -1131.5833
$ws.Range("H51").Value = 1625
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = -1016
$ws.Range("N51").Value = -2968
$ws.Range("H55").Value = 877.1429000000001
$ws.Range("I55").Value = 1605.7142
$ws.Range("J55").Value = 148.57143
$ws.Range("K55").Value = 1605.7142
$ws.Range("L55").Value = 148.57143
$ws.Range("M55").Value = -1391.7142
$ws.Range("N55").Value = -576.57143
$ws.Range("H69").Value = 4561.8335
$ws.Range("I69").Value = 4508.6665
$ws.Range("K69").Value = 13525.9995
$ws.Range("M69").Value = -12651.9995
$ws.Range("H72").Value = 4561.8335
$ws.Range("I72").Value = 4508.6665
$ws.Range("K72").Value = 40577.9985
$ws.Range("M72").Value = -36209.9985
$ws.Range("H100").Value = 4037.5
$ws.Range("I100").Value = 3337.5
$ws.Range("K100").Value = 3337.5
$ws.Range("M100").Value = -2796.5
$ws.Range("H111").Value = 868.5714
$ws.Range("I111").Value = 795.8
$ws.Range("K111").Value = 2387.4
$ws.Range("M111").Value = 679.6000000000004
$ws.Range("H113").Value = 4424.4443
$ws.Range("I113").Value = 4453.3335
$ws.Range("J113").Value = 4366.6665
$ws.Range("K113").Value = 4453.3335
$ws.Range("L113").Value = 4366.6665
$ws.Range("M113").Value = -1199.3335
$ws.Range("N113").Value = -10874.6665
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 26620
$ws.Range("J134").Value = 26620
$ws.Range("L134").Value = 26620
$ws.Range("N134").Value = -36760
$ws.Range("H136").Value = 29666.666
$ws.Range("J136").Value = 29666.666
$ws.Range("L136").Value = 29666.666
$ws.Range("N136").Value = -39866.666
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 17891.666
$ws.Range("J141").Value = 17891.666
$ws.Range("L141").Value = 17891.666
$ws.Range("N141").Value = -28251.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2862.375
$ws.Range("H131").Value = 1064.1897
$ws.Range("J131").Value = 1113.6415
$ws.Range("L131").Value = 3340.9245
$ws.Range("N131").Value = -13420.9245
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 58500
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -63560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3831.5833
$ws.Range("I80").Value = 3364.3333
$ws.Range("K80").Value = 3364.3333
$ws.Range("M80").Value = -2366.3333
$ws.Range("H83").Value = 3831.5833
$ws.Range("I83").Value = 3364.3333
$ws.Range("K83").Value = 16821.6665
$ws.Range("M83").Value = -11829.6665
$ws.Range("H97").Value = 1811.8
$ws.Range("I97").Value = 874.2857
$ws.Range("J97").Value = 3999.3333
$ws.Range("K97").Value = 874.2857
$ws.Range("L97").Value = 3999.3333
$ws.Range("M97").Value = -378.2857
$ws.Range("N97").Value = -4991.3333
$ws.Range("H133").Value = 29750
$ws.Range("J133").Value = 29750
$ws.Range("L133").Value = 29750
$ws.Range("N133").Value = -39870
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3227060
$ws.Range("I7").Value = 5263917
$ws.Range("K7").Value = 5263917
$ws.Range("M7").Value = -5263805
$ws.Range("H82").Value = 3236.7273
$ws.Range("I82").Value = 2451
$ws.Range("J82").Value = 3685.7144
$ws.Range("K82").Value = 2451
$ws.Range("L82").Value = 3685.7144
$ws.Range("M82").Value = -2090
$ws.Range("N82").Value = -4407.7144
$ws.Range("H85").Value = 3236.7273
$ws.Range("I85").Value = 2451
$ws.Range("J85").Value = 3685.7144
$ws.Range("K85").Value = 2451
$ws.Range("L85").Value = 3685.7144
$ws.Range("M85").Value = -1203
$ws.Range("N85").Value = -6181.7144
$ws.Range("H101").Value = 17454
$ws.Range("J101").Value = 17454
$ws.Range("L101").Value = 17454
$ws.Range("N101").Value = -23944
$ws.Range("H126").Value = 3227060
$ws.Range("I126").Value = 5263917
$ws.Range("K126").Value = 15791751
$ws.Range("M126").Value = -15789281

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2193.75
$ws.Range("J96").Value = 2430
$ws.Range("L96").Value = 2430
$ws.Range("N96").Value = -5176
$ws.Range("H100").Value = 576.8889
$ws.Range("I100").Value = 530.6667
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 1061.3334
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -520.3334
$ws.Range("N100").Value = -2282
$ws.Range("H103").Value = 25747.578
$ws.Range("J103").Value = 25747.578
$ws.Range("L103").Value = 25747.578
$ws.Range("N103").Value = -28091.578
$ws.Range("H132").Value = 139616.58
$ws.Range("I132").Value = 167463.45
$ws.Range("K132").Value = 502390.35
$ws.Range("M132").Value = -499860.35
$ws.Range("H135").Value = 150715
$ws.Range("J135").Value = 150715
$ws.Range("L135").Value = 150715
$ws.Range("N135").Value = -160855
$ws.Range("H137").Value = 37080
$ws.Range("J137").Value = 37080
$ws.Range("L137").Value = 37080
$ws.Range("N137").Value = -47280
$ws.Range("H139").Value = 29700
$ws.Range("J139").Value = 29700
$ws.Range("L139").Value = 29700
$ws.Range("N139").Value = -39980
